# digicode.xlsx — "Add files via upload" commit replay
#
# Adds one new product row (Microsoft 365 ESD licensing, OneDrive image) to
# the bottom of the "Tabla1" table on Sheet1, tweaks a couple of existing
# PRECIO EN SOLES / STOCK values on rows 7-8, updates the conditional
# formatting range to cover the new row, and leaves the selection sitting on
# C8 (no frozen/scrolled topLeftCell) the way the author's Excel session did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Small value edits on existing rows (rows 7 & 8)
# ---------------------------------------------------------------------
$ws.Range("H7").Value = 10
$ws.Range("C8").Value = 10
$ws.Range("H8").Value = 100

# ---------------------------------------------------------------------
# 2) Append a new row to the Tabla1 table (grows A1:J74 -> A1:J75)
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Tabla1")
$newListRow = $lo.ListRows.Add()
$newRow = $newListRow.Range.Row

# Pull formats from the row directly above (row 74) so the new row matches
# the rest of the table visually, then fill in the real values/formulas.
$ws.Range("A74").Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)
$ws.Range("B74").Copy()
$ws.Range("B" + $newRow).PasteSpecial(-4122)
$ws.Range("C74").Copy()
$ws.Range("C" + $newRow).PasteSpecial(-4122)
$ws.Range("E74").Copy()
$ws.Range("E" + $newRow).PasteSpecial(-4122)
$ws.Range("F74").Copy()
$ws.Range("F" + $newRow).PasteSpecial(-4122)
$ws.Range("G74").Copy()
$ws.Range("G" + $newRow).PasteSpecial(-4122)
$ws.Range("H74").Copy()
$ws.Range("H" + $newRow).PasteSpecial(-4122)
$ws.Range("I74").Copy()
$ws.Range("I" + $newRow).PasteSpecial(-4122)
$ws.Range("J74").Copy()
$ws.Range("J" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Calculated columns (CODIGO / MINI CODIGO)
$ws.Range("A" + $newRow).Formula = '="DIG"&TEXT(ROW()-1+100,"000")'
$ws.Range("F" + $newRow).Formula = "=+Tabla1[[#This Row],[CODIGO]]"

# Product data
$ws.Range("B" + $newRow).Value = "Licenciamiento Virtual (ESD) Microsoft 365 - 12 meses, 5 Dispositivos, 1TB de almacenamiento "
$ws.Range("C" + $newRow).Value = 10
$ws.Range("E" + $newRow).Value = "PRODUCTIVIDAD"
$ws.Range("G" + $newRow).Value = "Sucripcion x 365 dias"
$ws.Range("H" + $newRow).Value = 295
$ws.Range("I" + $newRow).Value = "https://licenciascol.com/cdn/shop/files/OneDrive.jpg"

# ---------------------------------------------------------------------
# 3) Extend the "STOCK < 1" conditional formatting down to the new row
# ---------------------------------------------------------------------
$oldCf = $ws.Range("C2:C74").FormatConditions.Item(1)
$oldCf.Delete()
$newCf = $ws.Range("C2:C" + $newRow).FormatConditions.Add(1, 6, "=1")
$newCf.Font.Color = -16383844
$newCf.Interior.Color = 13551615

# ---------------------------------------------------------------------
# 4) Selection / view state left by the author on save
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("C8").Select()
